$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data (rows 5-8), matching existing table layout:
# A: Placas, B: Fecha, C: Tipo de Gasto, D: Descripcion, E: Monto, F: Estado, G: Tiempo de Reparacion, H: Hora de Agregado
# An empty string "" marker means the source cell is blank (maps to the shared empty string).
$rows = @(
    @("ORB234", "17/12/2024", "Reparación",      "", 100.0, "FUNCIONAL", "", "03:37:34"),
    @("ORB234", "17/12/2024", "Reparación",      "", 22.0,  "FUNCIONAL", "", "03:38:26"),
    @("ORB234", "17/12/2024", "No Especificado", "u", 45.0, "FUNCIONAL", "", "03:38:55"),
    @("ORB234", "31/12/2024", "Mantenimiento",   "", 222.0, "FUNCIONAL", "", "03:44:57")
)

$startRow = 5
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]

    if ($rowData[3] -eq "") {
        $ws.Cells.Item($r, 4).Value = "'"
        $ws.Cells.Item($r, 4).ClearFormats()
    } else {
        $ws.Cells.Item($r, 4).Value = $rowData[3]
    }

    $ws.Cells.Item($r, 5).Value = $rowData[4]
    $ws.Cells.Item($r, 6).Value = $rowData[5]

    if ($rowData[6] -eq "") {
        $ws.Cells.Item($r, 7).Value = "'"
        $ws.Cells.Item($r, 7).ClearFormats()
    } else {
        $ws.Cells.Item($r, 7).Value = $rowData[6]
    }

    $ws.Cells.Item($r, 8).Value = $rowData[7]
}
